$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 237.03703
$ws.Range("I41").Value = 224.36363
$ws.Range("K41").Value = 224.36363
$ws.Range("M41").Value = 215.63637
# Row 76
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 113
$ws.Range("H113").Value = 8075.273
$ws.Range("I113").Value = 8242.125
$ws.Range("J113").Value = 7630.3335
$ws.Range("K113").Value = 8242.125
$ws.Range("L113").Value = 7630.3335
$ws.Range("M113").Value = -4988.125
$ws.Range("N113").Value = -14138.3335
# Row 132
$ws.Range("H132").Value = 910020.8
$ws.Range("I132").Value = 1022.9
$ws.Range("K132").Value = 3068.7
$ws.Range("M132").Value = -538.6999999999998
# Row 138
$ws.Range("H138").Value = 3741.93
$ws.Range("J138").Value = 3868.8525
$ws.Range("L138").Value = 11606.5575
$ws.Range("N138").Value = -21886.5575

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16663.072
$ws.Range("J32").Value = 59861.11
$ws.Range("L32").Value = 59861.11
$ws.Range("N32").Value = -60435.11
# Row 61
$ws.Range("H61").Value = 5340.88
$ws.Range("I61").Value = 4377.875
$ws.Range("K61").Value = 4377.875
$ws.Range("M61").Value = -4165.875
# Row 74
$ws.Range("H74").Value = 1663
$ws.Range("I74").Value = 1663
$ws.Range("K74").Value = 1663
$ws.Range("M74").Value = -789
# Row 77
$ws.Range("H77").Value = 1663
$ws.Range("I77").Value = 1663
$ws.Range("K77").Value = 8315
$ws.Range("M77").Value = -3947
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
# Row 132
$ws.Range("H132").Value = 2301
$ws.Range("I132").Value = 2009.711
$ws.Range("K132").Value = 6029.133
$ws.Range("M132").Value = -3499.133
# Row 136
$ws.Range("H136").Value = 5340.88
$ws.Range("I136").Value = 4377.875
$ws.Range("K136").Value = 13133.625
$ws.Range("M136").Value = -10583.625

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3692.3845
$ws.Range("I105").Value = 3106.1428
$ws.Range("J105").Value = 4376.3335
$ws.Range("K105").Value = 3106.1428
$ws.Range("L105").Value = 4376.3335
$ws.Range("M105").Value = -1359.1428
$ws.Range("N105").Value = -7870.3335
# Row 134
$ws.Range("H134").Value = 5827.4346
$ws.Range("I134").Value = 4672.231
$ws.Range("K134").Value = 14016.693
$ws.Range("M134").Value = -11481.693

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4768
$ws.Range("J16").Value = 3199.3333
$ws.Range("L16").Value = 3199.3333
$ws.Range("N16").Value = -3773.3333
# Row 113
$ws.Range("H113").Value = 4768
$ws.Range("J113").Value = 3199.3333
$ws.Range("L113").Value = 3199.3333
$ws.Range("N113").Value = -7539.3333
# Row 132
$ws.Range("H132").Value = 4155.75
$ws.Range("I132").Value = 3766.1333
$ws.Range("K132").Value = 11298.3999
$ws.Range("M132").Value = -8768.3999
# Row 134
$ws.Range("H134").Value = 2053.6
$ws.Range("J134").Value = 2170.8572
$ws.Range("L134").Value = 6512.571599999999
$ws.Range("N134").Value = -11582.5716

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 300
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -731
$ws.Range("N17").ClearContents()
# Row 46
$ws.Range("H46").Value = 1559.579
$ws.Range("I46").Value = 357
$ws.Range("J46").Value = 3621.1428
$ws.Range("K46").Value = 1071
$ws.Range("L46").Value = 10863.4284
$ws.Range("M46").Value = -980
$ws.Range("N46").Value = -11045.4284
# Row 107
$ws.Range("H107").Value = 941.6111
$ws.Range("J107").Value = 941.6111
$ws.Range("L107").Value = 2824.8333
$ws.Range("N107").Value = -6664.8333
# Row 113
$ws.Range("H113").Value = 1088.0588
$ws.Range("I113").Value = 1032.3334
$ws.Range("K113").Value = 3097.0002
$ws.Range("M113").Value = -927.0001999999999
# Row 122
$ws.Range("H122").Value = 1379.5454
$ws.Range("J122").Value = 1561.3636
$ws.Range("L122").Value = 14052.2724
$ws.Range("N122").Value = -18952.2724
# Row 137
$ws.Range("H137").Value = 3703.1667
$ws.Range("I137").Value = 3506.8572
$ws.Range("K137").Value = 10520.5716
$ws.Range("M137").Value = -5420.571599999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 32839.207
$ws.Range("I102").Value = 38994.418
$ws.Range("K102").Value = 38994.418
$ws.Range("M102").Value = -37372.418
# Row 113
$ws.Range("H113").Value = 3300
$ws.Range("I113").Value = 3166.6667
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3166.6667
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -996.6667000000002
$ws.Range("N113").Value = -7840
# Row 122
$ws.Range("H122").Value = 3429.1765
$ws.Range("I122").Value = 3098.6667
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 9296.000100000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6846.000100000001
$ws.Range("N122").Value = -15400
# Row 132
$ws.Range("H132").Value = 6389
$ws.Range("I132").Value = 6389
$ws.Range("K132").Value = 19167
$ws.Range("M132").Value = -16637

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9439.682000000001
$ws.Range("J7").Value = 9856.071
$ws.Range("L7").Value = 9856.071
$ws.Range("N7").Value = -10080.071
# Row 16
$ws.Range("H16").Value = 2362.1428
$ws.Range("I16").Value = 2362.1428
$ws.Range("K16").Value = 2362.1428
$ws.Range("M16").Value = -2192.1428
# Row 55
$ws.Range("H55").Value = 723.9091
$ws.Range("I55").Value = 737.4
$ws.Range("J55").Value = 589
$ws.Range("K55").Value = 737.4
$ws.Range("L55").Value = 589
$ws.Range("M55").Value = -564.4
$ws.Range("N55").Value = -935
# Row 61
$ws.Range("H61").Value = 9999.25
$ws.Range("J61").Value = 10004.5
$ws.Range("L61").Value = 10004.5
$ws.Range("N61").Value = -10408.5
# Row 68
$ws.Range("H68").Value = 991.3333
$ws.Range("I68").Value = 991.3333
$ws.Range("K68").Value = 991.3333
$ws.Range("M68").Value = -242.3333
# Row 71
$ws.Range("H71").Value = 991.3333
$ws.Range("I71").Value = 991.3333
$ws.Range("K71").Value = 4956.6665
$ws.Range("M71").Value = -1212.6665
# Row 82
$ws.Range("H82").Value = 1246.5834
$ws.Range("I82").Value = 936.1667
$ws.Range("J82").Value = 1557
$ws.Range("K82").Value = 936.1667
$ws.Range("L82").Value = 1557
$ws.Range("M82").Value = -575.1667
$ws.Range("N82").Value = -2279
# Row 85
$ws.Range("H85").Value = 1246.5834
$ws.Range("I85").Value = 936.1667
$ws.Range("J85").Value = 1557
$ws.Range("K85").Value = 936.1667
$ws.Range("L85").Value = 1557
$ws.Range("M85").Value = 311.8333
$ws.Range("N85").Value = -4053
# Row 113
$ws.Range("H113").Value = 9999.25
$ws.Range("J113").Value = 10004.5
$ws.Range("L113").Value = 10004.5
$ws.Range("N113").Value = -14344.5
# Row 126
$ws.Range("H126").Value = 9439.682000000001
$ws.Range("J126").Value = 9856.071
$ws.Range("L126").Value = 29568.213
$ws.Range("N126").Value = -34508.213
# Row 136
$ws.Range("H136").Value = 3458.25
$ws.Range("I136").Value = 3437.6843
$ws.Range("K136").Value = 10313.0529
$ws.Range("M136").Value = -7763.052899999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 37
$ws.Range("H37").Value = 21856.428
$ws.Range("I37").Value = 26999.5
$ws.Range("K37").Value = 26999.5
$ws.Range("M37").Value = -26796.5
# Row 40
$ws.Range("H40").Value = 39500
$ws.Range("I40").Value = 49000
$ws.Range("K40").Value = 49000
$ws.Range("M40").Value = -48851
# Row 42
$ws.Range("H42").Value = 49998
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
# Row 122
$ws.Range("H122").Value = 41487.832
$ws.Range("I122").Value = 41487.832
$ws.Range("K122").Value = 124463.496
$ws.Range("M122").Value = -122013.496
# Row 126
$ws.Range("H126").Value = 1858.8636
$ws.Range("J126").Value = 2511
$ws.Range("L126").Value = 7533
$ws.Range("N126").Value = -12473
# Row 136
$ws.Range("H136").Value = 2470.9143
$ws.Range("I136").Value = 2057.8635
$ws.Range("K136").Value = 6173.5905
$ws.Range("M136").Value = -3623.5905
